$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Colombian Primera B) keeps its League/Date/Time/Home/Away text,
#     but all the odds columns (F:AO, i.e. columns 6..41) get new values. ---
$row2 = @(1.01, 990, 19.5, 990, 1.2, 1.31, 0, 3.95, 1.01, 42, 1.01, 95, 1.01, 1.08, 1.68, 1.01, 1.01, 1.01, 1.32, 850, 1000, 1000, 850, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 6 + $i).Value = $row2[$i]
}

# --- Row 3 becomes the Colombian Primera A match (previously on row 4),
#     with a fresh set of odds, replacing the old Brazilian Serie A row.
#     B3 already holds the correct date text ("2025-11-12") and is left
#     untouched so it doesn't get reinterpreted as a date serial. ---
$ws.Range("A3").Value = "Colombian Primera A"
$ws.Range("C3").Value = "22:20:00"
$ws.Range("D3").Value = "Boyaca Chico"
$ws.Range("E3").Value = "Millonarios"

$row3 = @(5.4, 5.9, 1.23, 1.77, 3.9, 4.3, 1.01, 1.07, 1.85, 1.01, 1.85, 1.01, 1.19, 1.01, 1.06, 2.04, 1.01, 1.01, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 42, 1000, 1000, 85, 260, 1000, 30)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 6 + $i).Value = $row3[$i]
}

# --- The old row 4 (previous Colombian Primera A entry) is removed entirely
#     now that its data lives on row 3; the sheet shrinks to 3 rows. ---
$ws.Rows(4).Delete()
